$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2022-09-10 20:58:28"

for ($row = 2; $row -le 64; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
